# Generate Report for Handoff
#
# A new handoff run produced a fresh GUID-named markdown source and fresh
# xliff target files; refresh the report's Overview / zh-cn / de-de sheets
# with the new file names and the new handoff timestamps. The hyperlinks
# that point at the source markdown file keep their original target URL
# (only their visible/display text tracks the new file name, matching the
# upstream report generator).

$wb = $excel.ActiveWorkbook

$newGuid   = "18d60560-cda4-4d1e-bebc-f28c2cfcc9f6"
$newZhHash = "e1bc6fc5843a46703f51fedfd218eff0c5ade03f"

# Hyperlink target (left untouched - it still points at the old commit's
# markdown file; only the display text for the new file name changes).
$hyperlinkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a96eedefdf729c91355cdd617b269f4242cb0f82/e2e/c2483ccd-5a9d-4f33-aa89-074a63a27c56.md"

function Update-SourceHyperlink($ws, [string]$cellRef, [string]$displayText) {
    # Refresh a hyperlink's cell value + display text in place without
    # leaving a stale duplicate hyperlink entry behind.
    $ws.Hyperlinks.Delete()
    $rng = $ws.Range($cellRef)
    $rng.Value = $displayText
    $ws.Hyperlinks.Add($rng, $hyperlinkAddress, "", "", $displayText) | Out-Null
}

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"
Update-SourceHyperlink $wsOverview "B2" "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-08-28 16:58:13"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
Update-SourceHyperlink $wsZhCn "A2" "$newGuid.md"
$wsZhCn.Range("G2").Value = "$newGuid.$newZhHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-28 16:58:08"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
Update-SourceHyperlink $wsDeDe "A2" "$newGuid.md"
$wsDeDe.Range("G2").Value = "$newGuid.$newZhHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-28 16:58:13"
